# Add "Scenario 4: Carry (Loss)" to the "Scenarios" worksheet, mirroring the
# existing "Scenario 3: Carry" block (rows 14-20) into new rows 21-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

# 1) Clone the formatting of the Scenario 3 block (blank separator row,
#    title row, header row, 4 data rows) onto the new rows so the new block
#    uses the exact same cell styles (s="5" / s="8" / s="9") already present
#    in the workbook. Copy row-by-row, matching the exact column footprint
#    of each source row (the title row only spans A:B).
$ws.Range("A14:G14").Copy()
$ws.Range("A21:G21").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A15:B15").Copy()
$ws.Range("A22:B22").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A16:G16").Copy()
$ws.Range("A23:G23").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A17:G20").Copy()
$ws.Range("A24:G27").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# Row 21 is left blank (separator row).
"A21","B21","C21","D21","E21","F21","G21" | ForEach-Object { $ws.Range($_).Value = "" }

# 2) Row 22: scenario title + description
$ws.Range("A22").Value = "Scenario 4: Carry (Loss)"
$ws.Range("B22").Value = "T1: P1A(1400) & P1B(600) vs T2: Avg 1000. 0-2 Loss for T1."

# 3) Row 23: table header
$ws.Range("A23").Value = "Player"
$ws.Range("B23").Value = "Start ELO"
$ws.Range("C23").Value = "Games"
$ws.Range("D23").Value = "Opp. Avg"
$ws.Range("E23").Value = "Result"
$ws.Range("F23").Value = "Delta"
$ws.Range("G23").Value = "End ELO"

# 4) Rows 24-27: player data
$data = @(
    @("P1A (1400)", 1400, 50, 1000, "Loss", -11, 1389),
    @("P1B (600)",   600, 50, 1000, "Loss",  -4,  596),
    @("P2A (1000)", 1000, 50, 1000, "Win",     7, 1007),
    @("P2B (1000)", 1000, 50, 1000, "Win",     7, 1007)
)

$cols = @("A","B","C","D","E","F","G")
$rowNum = 24
foreach ($row in $data) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$rowNum").Value = $row[$i]
    }
    $rowNum++
}
